$wb = $excel.ActiveWorkbook

# --- Add the new "Terms" worksheet right after Sheet1 ---
$sheet1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add([Type]::Missing, $sheet1)
$ws.Name = "Terms"

# --- Fill in the license / attribution text ---
$ws.Range("A1").Value = "This dataset on 'SleepData.xlsx' is hypothetical and was generated"
$ws.Range("A2").Value = "by Paolo G. Hilado (Github: Dcroix) for training purposes on Basic Statistics . Considering"
$ws.Range("A3").Value = "that most of the values generated by this dataset use randomization, "
$ws.Range("A4").Value = "in such a rare case that it resembles any existing dataset, it is purely "
$ws.Range("A5").Value = "coincidental. It is distributed under "
$ws.Range("A6").Value = " Creative Commons Attribution-NoDerivatives 4.0 International Public License."

# --- Column A is wide, to hold the long lines of text ---
$ws.Columns.Item(1).ColumnWidth = 84.6667

# --- Turn the last line into a hyperlink to the license on GitHub ---
$licenseUrl = "https://github.com/Dcroix/SampleData/blob/master/Creative Commons Attribution-NoDerivatives 4.0 International Public License"
$ws.Hyperlinks.Add(
    $ws.Range("A6"),
    $licenseUrl,
    [Type]::Missing,
    [Type]::Missing,
    $licenseUrl
)
# Hyperlinks.Add() also overwrote the cell text with the display text we just
# passed it (the URL) - put the real sentence back now that the hyperlink
# (and its "display" attribute) have been recorded against the cell.
$ws.Range("A6").Value = " Creative Commons Attribution-NoDerivatives 4.0 International Public License."

# --- Match the saved selection/activation state from the source file ---
$ws.Range("A10").Select() | Out-Null
